$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B10").Value = "www.stat.gov.kg"
$ws.Range("B9").Value = "(0312) 62 56 07"
$ws.Range("B7").Value = "Mambetaliev T.A."
$ws.Range("B6").Value = "National Statistical Committee of the Kyrgyz Republic" + [char]10 + "Department of Digital Development and Sustainable Development Statistics"

$ws.Range("B9").Select()
